# Append 45 new rows (102-146) to the master-reg_center_machine_device sheet,
# continuing the existing regcntr_id/machine_id/device_id pattern, and leave
# the newly-added block selected (matching the author's last on-screen state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 102
$endRow   = 146
$startC   = 3000121

for ($i = 0; $i -lt ($endRow - $startRow + 1); $i++) {
    $row = $startRow + $i
    $cycle = $i % 9

    $ws.Cells.Item($row, 1).Value = 10002 + $cycle
    $ws.Cells.Item($row, 2).Value = 10021 + $cycle
    $ws.Cells.Item($row, 3).Value = $startC + $i
    $ws.Cells.Item($row, 4).Value = "eng"
    $ws.Cells.Item($row, 5).Value = $true
    $ws.Cells.Item($row, 6).Value = "superadmin()"
    $ws.Cells.Item($row, 7).Value = "now()"
}

$ws.Range("A$startRow`:G$endRow").Select()

# Page setup was touched to portrait orientation in the saved file as well.
$ws.PageSetup.Orientation = 1
